## Update "current calculation" values (D column) for rows 8-35, add the
## missing D/E/F data for rows 35-47, extend the shared E/F formulas, and
## move the active selection/scroll position — per the "current calculation
## done as per required" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Revised "acutal amp" (D column) readings for the existing rows 8-34
# ---------------------------------------------------------------------
$dUpdates = @{
    8  = 0.97
    9  = 1.24
    10 = 1.7
    11 = 2.06
    12 = 2.52
    13 = 2.9
    14 = 3.28
    15 = 3.85
    16 = 4.36
    17 = 4.88
    18 = 5.33
    19 = 5.83
    20 = 6.47
    21 = 6.92
    22 = 7.51
    23 = 8.01
    24 = 8.76
    25 = 9.27
    26 = 10
    27 = 10.42
    28 = 10.96
    29 = 11.51
    30 = 12.04
    31 = 12.69
    32 = 13.14
    33 = 13.8
    34 = 14.32
}

foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# ---------------------------------------------------------------------
# 2) Row 35 used to divide-by-zero (D35 was blank) - give it real data
#    so E35/F35 resolve instead of showing #DIV/0!
# ---------------------------------------------------------------------
$ws.Cells.Item(35, 4).Value = 14.96

# ---------------------------------------------------------------------
# 3) Rows 36-47 previously only had the "current meter" (B) value - fill
#    in the matching "acutal amp" (D) and restore the calc/multiply
#    formulas (E/F) that the rest of the table uses.
# ---------------------------------------------------------------------
$dNew = @{
    36 = 15.43
    37 = 16.04
    38 = 16.61
    39 = 17.09
    40 = 17.69
    41 = 18.17
    42 = 18.78
    43 = 19.29
    44 = 19.79
    45 = 20.63
    46 = 21.23
    47 = 21.92
}

foreach ($row in $dNew.Keys) {
    $ws.Cells.Item($row, 4).Value = $dNew[$row]
}

# Fill E36:E45 / F36:F45 as one pass (mirrors the existing E7:E35 /
# F7:F35 "calulation value" + "make actual value with multiply" formulas)
$ws.Range("E36:E45").Formula = "=B36/D36"
$ws.Range("F36:F45").Formula = "=D36*E36"

# Rows 46:47 were filled in a separate pass, so they pick up their own
# shared-formula group, same as the source workbook.
$ws.Range("E46:E47").Formula = "=B46/D46"
$ws.Range("F46:F47").Formula = "=D46*E46"

# ---------------------------------------------------------------------
# 4) Move the saved view: scroll down to row 25 and land the selection
#    on E7 (was A4 / D26).
# ---------------------------------------------------------------------
$ws.Range("E7").Select()
try {
    $excel.ActiveWindow.ScrollRow = 25
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$wb.Application.Calculate()
